# Update the "想去人数" (want-to-go count) values on the 展览 and 全部类型
# sheets, matching the data refresh reflected in the commit
# "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value  = 14847
$wsExhibit.Range("F3").Value  = 18427
$wsExhibit.Range("F14").Value = 99
$wsExhibit.Range("F15").Value = 196
$wsExhibit.Range("F17").Value = 1403
$wsExhibit.Range("F22").Value = 7634
$wsExhibit.Range("F23").Value = 986
$wsExhibit.Range("F28").Value = 5943
$wsExhibit.Range("F34").Value = 5279

# Sheet "全部类型" (All types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value  = 14847
$wsAll.Range("F3").Value  = 18427
$wsAll.Range("F14").Value = 99
$wsAll.Range("F15").Value = 196
$wsAll.Range("F17").Value = 1403
$wsAll.Range("F23").Value = 7634
$wsAll.Range("F24").Value = 986
$wsAll.Range("F30").Value = 5943
$wsAll.Range("F36").Value = 5279
